# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 51 (pushing the existing
# "Terminal La Palmera de La Serena - Albahaca" records down by one row),
# then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 51; rows 51:112 shift down to 52:113.
$ws.Rows("51:51").Insert()

# Populate the new row 51 with the new weekly record.
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44799
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = 100112052
$ws.Cells.Item(51, 7).Value = "Albahaca"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 1000
$ws.Cells.Item(51, 11).Value = 4000
$ws.Cells.Item(51, 12).Value = 4500
$ws.Cells.Item(51, 13).Value = 4250
$ws.Cells.Item(51, 14).Value = "$/paquete"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 4250
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
